$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.69037401676178
$ws.Range("B1").Value = 3.327393054962158
$ws.Range("C1").Value = 6.210051536560059
$ws.Range("D1").Value = 1.802961826324463
$ws.Range("E1").Value = 0.8911571502685547
